$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "46.184.90"
$ws.Range("E2").Value = "  -0.89%  "

# Row 3
$ws.Range("D3").Value = "2.609.37"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "

# Row 7
$ws.Range("E7").Value = "  -0.69%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("E9").Value = "  +1.54%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.95%  "

# Row 12
$ws.Range("E12").Value = "  +0.41%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.17%  "

# Row 14
$ws.Range("D14").Value = "3.005.22"
$ws.Range("E14").Value = "  +0.87%  "

# Row 15
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
$ws.Range("D16").Value = "2.601.91"
$ws.Range("E16").Value = "  +0.69%  "

# Row 17
$ws.Range("E17").Value = "  +1.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("D19").Value = "46.344.12"
$ws.Range("E19").Value = "  -0.89%  "

# Row 20
$ws.Range("E20").Value = "  +0.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.53%  "

# Row 22
$ws.Range("E22").Value = "  +1.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "276.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.72%  "

# Row 26
$ws.Range("E26").Value = "  +0.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "

# Row 29
$ws.Range("E29").Value = "  +1.25%  "

# Row 30
$ws.Range("E30").Value = "  +2.72%  "

# Row 31 - coin swapped with row 32 (InjectiveProtocol -> Toncoin)
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.33%  "

# Row 32 - coin swapped with row 31 (Toncoin -> InjectiveProtocol)
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.47%  "

# Row 33
$ws.Range("E33").Value = "  +1.92%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.91%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0840"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.82%  "

# Row 39
$ws.Range("E39").Value = "  +5.57%  "

# Row 40
$ws.Range("E40").Value = "  +0.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.11%  "

# Row 43
$ws.Range("E43").Value = "  +2.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.98%  "

# Row 46
$ws.Range("D46").Value = "2.106.66"
$ws.Range("E46").Value = "  +4.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.58%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "109.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "

# Row 51 - Stacks replaced with Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.201"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
